$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column (Price) cells: force text format so numeric-looking strings
# (e.g. "1.00", "23.40", "0.0168") are not coerced into numbers, which
# would silently change/lose their displayed text. Reset the style back
# to Normal afterwards so only the cell VALUE changes (matching the diff).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.631.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.635.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("E6").Value = "  -0.47%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.40"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.86%  "

$ws.Range("E9").Value = "  +2.70%  "

$ws.Range("E10").Value = "  +0.31%  "

$ws.Range("E11").Value = "  -2.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.866.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.635.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.617.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.21%  "

$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("E20").Value = "  -1.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.44%  "

$ws.Range("E23").Value = "  +1.75%  "

$ws.Range("E24").Value = "  +7.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.62%  "

$ws.Range("E26").Value = "  -0.66%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("E32").Value = "  -0.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.478.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.44%  "

$ws.Range("E36").Value = "  -1.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.943"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.884"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.92%  "

$ws.Range("E39").Value = "  -0.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0168"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.47%  "

$ws.Range("E41").Value = "  +2.00%  "

$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("E43").Value = "  -2.28%  "

$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("E45").Value = "  -0.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.776.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("E48").Value = "  +2.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.62%  "

$ws.Range("E50").Value = "  -1.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0995"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.11%  "
